$d = $word.ActiveDocument

$d.Paragraphs.Item(2).Range.Text = "Texte 2"
$d.Paragraphs.Item(3).Range.Text = "Texte 3"
